# musicas.xlsx — row 536 ("Jazz00054") had been removed earlier, leaving an
# empty placeholder row (only A536/D536 present, no value). This edit closes
# that gap: every data row from 537 down to 681 is shifted up by one in the
# title/author/genre columns (B:D) while the running id column (A) is left
# alone since it was already sequential. The now-unused last row (681) is
# cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstGapRow = 536
$lastDataRow = 681

# Shift columns B (titulo), C (autor) and D (genero) up by one row, starting
# at the gap and walking down to the last populated row. Processing top to
# bottom is safe because each row's new content is read from row+1 *before*
# that row has been overwritten.
for ($r = $firstGapRow; $r -lt $lastDataRow; $r++) {
    $titulo = $ws.Cells.Item($r + 1, 2).Value2
    $autor  = $ws.Cells.Item($r + 1, 3).Value2
    $genero = $ws.Cells.Item($r + 1, 4).Value2

    $ws.Cells.Item($r, 2).Value = $titulo
    $ws.Cells.Item($r, 3).Value = $autor
    $ws.Cells.Item($r, 4).Value = $genero
}

# The gap row's id cell already existed (empty); fill it in now that it holds
# real data (id = row - 1).
$ws.Cells.Item($firstGapRow, 1).Value = $firstGapRow - 1

# The final row has nothing left to show after the shift: drop its id/genre
# cells entirely and blank out title/author back to the sheet's default
# (un-styled) look.
$ws.Cells.Item($lastDataRow, 1).Clear()
$ws.Cells.Item($lastDataRow, 4).Clear()

$ws.Cells.Item($lastDataRow, 2).Clear()
$ws.Cells.Item($lastDataRow, 3).Clear()
$ws.Cells.Item($lastDataRow, 2).Value = ""
$ws.Cells.Item($lastDataRow, 3).Value = ""
$ws.Cells.Item($lastDataRow, 2).NumberFormat = "General"
$ws.Cells.Item($lastDataRow, 3).NumberFormat = "General"

# Restore the view: selection sits on C536 (the newly-filled cell) and the
# window is scrolled up a few rows from where it was.
$ws.Range("C536").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 525
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "musicas.xlsx: closed gap at row $firstGapRow, cleared row $lastDataRow"
